# FA_IMPORT worksheet cleanup: remove two duplicate rows from the
# PRIMARY_TERM_ADDRESS / FAC_ADDR list (row 208 = "105 HICKORY LN",
# row 192 = "1302 WALKER DR"), then bring the AutoFilter range, the
# hidden _FilterDatabase defined name, and the active selection back
# in sync with the new (smaller) used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the lower-numbered duplicate row first so the row number of
# the other deletion target isn't shifted before it is removed.
$ws.Rows(208).Delete()
$ws.Rows(192).Delete()

# Re-apply AutoFilter so its stored range shrinks from A1:E645 to A1:E643.
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:E643").AutoFilter()

# Keep the hidden _FilterDatabase name (driven by the AutoFilter above)
# in sync with the new used range as well.
$wb.Names.Item("FA_IMPORT!_FilterDatabase").RefersTo = "=FA_IMPORT!`$A`$1:`$E`$643"

# Restore the active selection to A202, matching the saved view state.
[void]$ws.Range("A202").Select()
